$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Business Segment" / "Import Action" / "Template Name" / "Section Name"
# block (rows 16-19), shifting the Product Group section (and everything below it)
# up by four rows.
$ws.Rows("16:19").Delete()

# The "WAC" price cells (now on rows 23 and 27 after the shift) were stored as
# text labels ("$195" / "$425"); make them real currency numbers with a
# currency number format instead.
$ws.Range("E23").Value = 195
$ws.Range("E23").NumberFormat = '"$"#,##0_);[Red]\("$"#,##0\)'

$ws.Range("E27").Value = 425
$ws.Range("E27").NumberFormat = '"$"#,##0_);[Red]\("$"#,##0\)'

# Match the saved selection / active cell.
[void]$ws.Range("G12").Select()

# Printed page orientation was set explicitly (portrait).
$ws.PageSetup.Orientation = 1
